$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 93126.63
$ws.Cells.Item(6, 9).Value = 93126.63
$ws.Cells.Item(6, 11).Value = 279379.89
$ws.Cells.Item(6, 13).Value = -279267.89
$ws.Cells.Item(17, 8).Value = 305376.94
$ws.Cells.Item(17, 10).Value = 347277.28
$ws.Cells.Item(17, 12).Value = 1041831.84
$ws.Cells.Item(17, 14).Value = -1042167.84
$ws.Cells.Item(41, 8).Value = 1797.2
$ws.Cells.Item(41, 9).Value = 316.66666
$ws.Cells.Item(41, 10).Value = 2431.7144
$ws.Cells.Item(41, 11).Value = 316.66666
$ws.Cells.Item(41, 12).Value = 2431.7144
$ws.Cells.Item(41, 13).Value = 123.33334
$ws.Cells.Item(41, 14).Value = -3311.7144
$ws.Cells.Item(62, 8).Value = 35720260
$ws.Cells.Item(62, 9).Value = 62505000
$ws.Cells.Item(62, 10).Value = 7276.3335
$ws.Cells.Item(62, 11).Value = 62505000
$ws.Cells.Item(62, 12).Value = 7276.3335
$ws.Cells.Item(62, 13).Value = -62504376
$ws.Cells.Item(62, 14).Value = -8524.333500000001
$ws.Cells.Item(65, 8).Value = 35720260
$ws.Cells.Item(65, 9).Value = 62505000
$ws.Cells.Item(65, 10).Value = 7276.3335
$ws.Cells.Item(65, 11).Value = 312525000
$ws.Cells.Item(65, 12).Value = 36381.6675
$ws.Cells.Item(65, 13).Value = -312521880
$ws.Cells.Item(65, 14).Value = -42621.6675
$ws.Cells.Item(69, 8).Value = 9622.637000000001
$ws.Cells.Item(69, 9).Value = 9622.637000000001
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 28867.911
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -27993.911
$ws.Cells.Item(72, 8).Value = 9622.637000000001
$ws.Cells.Item(72, 9).Value = 9622.637000000001
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 86603.73300000001
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -82235.73300000001
$ws.Cells.Item(98, 8).Value = 3481377.8
$ws.Cells.Item(98, 9).Value = 3637681
$ws.Cells.Item(98, 11).Value = 3637681
$ws.Cells.Item(98, 13).Value = -3636183
$ws.Cells.Item(111, 8).Value = 3957.3572
$ws.Cells.Item(111, 10).Value = 5529.375
$ws.Cells.Item(111, 12).Value = 16588.125
$ws.Cells.Item(111, 14).Value = -22722.125
$ws.Cells.Item(115, 8).Value = 10351647
$ws.Cells.Item(115, 9).Value = 12075254
$ws.Cells.Item(115, 11).Value = 36225762
$ws.Cells.Item(115, 13).Value = -36224195
$ws.Cells.Item(116, 8).Value = 16212.857
$ws.Cells.Item(116, 9).Value = 19198.545
$ws.Cells.Item(116, 10).Value = 5265.3335
$ws.Cells.Item(116, 11).Value = 19198.545
$ws.Cells.Item(116, 12).Value = 5265.3335
$ws.Cells.Item(116, 13).Value = -15756.545
$ws.Cells.Item(116, 14).Value = -12149.3335
$ws.Cells.Item(122, 8).Value = 3481377.8
$ws.Cells.Item(122, 9).Value = 3637681
$ws.Cells.Item(122, 11).Value = 10913043
$ws.Cells.Item(122, 13).Value = -10910593
$ws.Cells.Item(138, 8).Value = 4865.97
$ws.Cells.Item(138, 9).Value = 1667.1177
$ws.Cells.Item(138, 10).Value = 5521.1567
$ws.Cells.Item(138, 11).Value = 5001.3531
$ws.Cells.Item(138, 12).Value = 16563.4701
$ws.Cells.Item(138, 13).Value = 138.6468999999997
$ws.Cells.Item(138, 14).Value = -26843.4701
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7070.493
$ws.Cells.Item(32, 9).Value = 7070.493
$ws.Cells.Item(32, 11).Value = 7070.493
$ws.Cells.Item(32, 13).Value = -6783.493
$ws.Cells.Item(61, 8).Value = 4155.4346
$ws.Cells.Item(61, 9).Value = 1960.2307
$ws.Cells.Item(61, 10).Value = 7009.2
$ws.Cells.Item(61, 11).Value = 1960.2307
$ws.Cells.Item(61, 12).Value = 7009.2
$ws.Cells.Item(61, 13).Value = -1748.2307
$ws.Cells.Item(61, 14).Value = -7433.2
$ws.Cells.Item(63, 8).Value = 2999.6667
$ws.Cells.Item(63, 9).Value = 2999
$ws.Cells.Item(63, 11).Value = 2999
$ws.Cells.Item(63, 13).Value = -2313
$ws.Cells.Item(66, 8).Value = 2999.6667
$ws.Cells.Item(66, 9).Value = 2999
$ws.Cells.Item(66, 11).Value = 14995
$ws.Cells.Item(66, 13).Value = -11563
$ws.Cells.Item(97, 8).Value = 797.05554
$ws.Cells.Item(97, 9).Value = 514.36365
$ws.Cells.Item(97, 10).Value = 1241.2858
$ws.Cells.Item(97, 11).Value = 514.36365
$ws.Cells.Item(97, 12).Value = 1241.2858
$ws.Cells.Item(97, 13).Value = -18.36365000000001
$ws.Cells.Item(97, 14).Value = -2233.2858
$ws.Cells.Item(110, 8).Value = 10050
$ws.Cells.Item(110, 9).Value = 8755.611000000001
$ws.Cells.Item(110, 11).Value = 8755.611000000001
$ws.Cells.Item(110, 13).Value = -6710.611000000001
$ws.Cells.Item(136, 8).Value = 4155.4346
$ws.Cells.Item(136, 9).Value = 1960.2307
$ws.Cells.Item(136, 10).Value = 7009.2
$ws.Cells.Item(136, 11).Value = 5880.6921
$ws.Cells.Item(136, 12).Value = 21027.6
$ws.Cells.Item(136, 13).Value = -3330.6921
$ws.Cells.Item(136, 14).Value = -26127.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).ClearContents()
$ws.Cells.Item(68, 14).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).ClearContents()
$ws.Cells.Item(71, 14).Value = 0
$ws.Cells.Item(105, 8).Value = 1959.091
$ws.Cells.Item(105, 9).Value = 1959.091
$ws.Cells.Item(105, 11).Value = 1959.091
$ws.Cells.Item(105, 13).Value = -212.0909999999999
$ws.Cells.Item(134, 8).Value = 1891.5103
$ws.Cells.Item(134, 9).Value = 1447.9524
$ws.Cells.Item(134, 11).Value = 4343.857199999999
$ws.Cells.Item(134, 13).Value = -1808.857199999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 11
$ws.Cells.Item(25, 9).Value = 11
$ws.Cells.Item(25, 11).Value = 11
$ws.Cells.Item(25, 13).Value = 163
$ws.Cells.Item(105, 8).Value = 1865.75
$ws.Cells.Item(105, 9).Value = 1477
$ws.Cells.Item(105, 10).Value = 4198.25
$ws.Cells.Item(105, 11).Value = 1477
$ws.Cells.Item(105, 12).Value = 4198.25
$ws.Cells.Item(105, 13).Value = 270
$ws.Cells.Item(105, 14).Value = -7692.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 3296
$ws.Cells.Item(16, 9).Value = 2999.5
$ws.Cells.Item(16, 11).Value = 8998.5
$ws.Cells.Item(16, 13).Value = -8825.5
$ws.Cells.Item(18, 8).Value = 745.0833
$ws.Cells.Item(18, 9).Value = 527.3333
$ws.Cells.Item(18, 11).Value = 1581.9999
$ws.Cells.Item(18, 13).Value = -1412.9999
$ws.Cells.Item(122, 8).Value = 1983.3334
$ws.Cells.Item(122, 10).Value = 1975
$ws.Cells.Item(122, 12).Value = 17775
$ws.Cells.Item(122, 14).Value = -22675
$ws.Cells.Item(132, 8).Value = 5854.04
$ws.Cells.Item(132, 9).Value = 6906.9
$ws.Cells.Item(132, 10).Value = 1642.6
$ws.Cells.Item(132, 11).Value = 62162.1
$ws.Cells.Item(132, 12).Value = 14783.4
$ws.Cells.Item(132, 13).Value = -59632.1
$ws.Cells.Item(132, 14).Value = -19843.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7376.4346
$ws.Cells.Item(70, 9).Value = 7892.125
$ws.Cells.Item(70, 11).Value = 7892.125
$ws.Cells.Item(70, 13).Value = -7622.125
$ws.Cells.Item(73, 8).Value = 7376.4346
$ws.Cells.Item(73, 9).Value = 7892.125
$ws.Cells.Item(73, 11).Value = 7892.125
$ws.Cells.Item(73, 13).Value = -6956.125
$ws.Cells.Item(80, 8).Value = 3482.5854
$ws.Cells.Item(80, 9).Value = 3619.0645
$ws.Cells.Item(80, 10).Value = 3059.5
$ws.Cells.Item(80, 11).Value = 3619.0645
$ws.Cells.Item(80, 12).Value = 3059.5
$ws.Cells.Item(80, 13).Value = -2621.0645
$ws.Cells.Item(80, 14).Value = -5055.5
$ws.Cells.Item(83, 8).Value = 3482.5854
$ws.Cells.Item(83, 9).Value = 3619.0645
$ws.Cells.Item(83, 10).Value = 3059.5
$ws.Cells.Item(83, 11).Value = 18095.3225
$ws.Cells.Item(83, 12).Value = 15297.5
$ws.Cells.Item(83, 13).Value = -13103.3225
$ws.Cells.Item(83, 14).Value = -25281.5
$ws.Cells.Item(112, 8).Value = 40293
$ws.Cells.Item(112, 10).Value = 40293
$ws.Cells.Item(112, 12).Value = 40293
$ws.Cells.Item(112, 14).Value = -42509
$ws.Cells.Item(122, 8).Value = 1806.3478
$ws.Cells.Item(122, 9).Value = 1803.1364
$ws.Cells.Item(122, 11).Value = 5409.4092
$ws.Cells.Item(122, 13).Value = -2959.4092
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 8253.223
$ws.Cells.Item(40, 9).Value = 7757.2856
$ws.Cells.Item(40, 11).Value = 7757.2856
$ws.Cells.Item(40, 13).Value = -7621.2856
$ws.Cells.Item(122, 8).Value = 7739.15
$ws.Cells.Item(122, 9).Value = 7714.846
$ws.Cells.Item(122, 10).Value = 7784.2856
$ws.Cells.Item(122, 11).Value = 23144.538
$ws.Cells.Item(122, 12).Value = 23352.8568
$ws.Cells.Item(122, 13).Value = -20694.538
$ws.Cells.Item(122, 14).Value = -28252.8568
$ws.Cells.Item(136, 8).Value = 4353.5
$ws.Cells.Item(136, 9).Value = 4245.25
$ws.Cells.Item(136, 10).Value = 4425.6665
$ws.Cells.Item(136, 11).Value = 12735.75
$ws.Cells.Item(136, 12).Value = 13276.9995
$ws.Cells.Item(136, 13).Value = -10185.75
$ws.Cells.Item(136, 14).Value = -18376.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 37501
$ws.Cells.Item(2, 9).Value = 30002
$ws.Cells.Item(2, 10).Value = 45000
$ws.Cells.Item(2, 11).Value = 30002
$ws.Cells.Item(2, 12).Value = 45000
$ws.Cells.Item(2, 13).Value = -29890
$ws.Cells.Item(2, 14).Value = -45224
$ws.Cells.Item(6, 8).Value = 1000
$ws.Cells.Item(6, 9).Value = 1000
$ws.Cells.Item(6, 11).Value = 1000
$ws.Cells.Item(6, 13).Value = -885
$ws.Cells.Item(51, 8).Value = 27993
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()
$ws.Cells.Item(126, 8).Value = 2475.8462
$ws.Cells.Item(126, 9).Value = 2380.7273
$ws.Cells.Item(126, 10).Value = 2999
$ws.Cells.Item(126, 11).Value = 7142.1819
$ws.Cells.Item(126, 12).Value = 8997
$ws.Cells.Item(126, 13).Value = -4672.1819
$ws.Cells.Item(126, 14).Value = -13937
$ws.Cells.Item(132, 8).Value = 1076.7894
$ws.Cells.Item(132, 9).Value = 1122.5834
$ws.Cells.Item(132, 11).Value = 3367.7502
$ws.Cells.Item(132, 13).Value = -837.7501999999999
